# Apply updated crypto price/volume values (cols D and E) to match latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.310.90"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "1.875.12"
$ws.Range("E3").Value = "  +0.02%  "
$ws.Range("E4").Value = "  +0.11%  "
$c = $ws.Range("D5")
$c.Value = "'0.7102"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.43%  "
$c = $ws.Range("D6")
$c.Value = "'241.74"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("E7").Value = "  +0.04%  "
$c = $ws.Range("D8")
$c.Value = "'0.07877"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +1.97%  "
$c = $ws.Range("D9")
$c.Value = "'0.3125"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +0.58%  "
$c = $ws.Range("D10")
$c.Value = "'25.22"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +0.60%  "
$c = $ws.Range("D11")
$c.Value = "'0.08404"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +0.22%  "
$ws.Range("D12").Value = "1.874.71"
$ws.Range("E12").Value = "  -0.29%  "
$c = $ws.Range("D13")
$c.Value = "'5.242"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +0.48%  "
$c = $ws.Range("D14")
$c.Value = "'0.7170"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +0.77%  "
$c = $ws.Range("D15")
$c.Value = "'91.22"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -0.14%  "
$c = $ws.Range("D16")
$c.Value = "'6.201"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +3.73%  "
$c = $ws.Range("D17")
$c.Value = "'0.000008358"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +0.97%  "
$ws.Range("D18").Value = "29.309.63"
$ws.Range("E18").Value = "  +0.03%  "
$c = $ws.Range("D19")
$c.Value = "'240.64"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -0.83%  "
$c = $ws.Range("D20")
$c.Value = "'13.22"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.16%  "
$ws.Range("D21").Value = "2.124.62"
$ws.Range("E21").Value = "  -0.18%  "
$c = $ws.Range("D22")
$c.Value = "'1.000"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.11%  "
$c = $ws.Range("D23")
$c.Value = "'7.794"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.31%  "
$c = $ws.Range("D24")
$c.Value = "'1.000"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.08%  "
$c = $ws.Range("D25")
$c.Value = "'0.1591"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -1.77%  "
$c = $ws.Range("D26")
$c.Value = "'162.66"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -0.37%  "
$c = $ws.Range("D27")
$c.Value = "'9.048"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +0.29%  "
$c = $ws.Range("D28")
$c.Value = "'18.52"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -0.02%  "
$c = $ws.Range("D29")
$c.Value = "'1.505"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +0.07%  "
$c = $ws.Range("D30")
$c.Value = "'4.422"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +0.00%  "
$c = $ws.Range("D31")
$c.Value = "'4.347"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +0.43%  "
$c = $ws.Range("D32")
$c.Value = "'1.204"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -6.62%  "
$c = $ws.Range("D33")
$c.Value = "'0.05356"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +2.10%  "
$c = $ws.Range("D34")
$c.Value = "'1.945"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +0.90%  "
$c = $ws.Range("D35")
$c.Value = "'0.7502"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +0.27%  "
$c = $ws.Range("D36")
$c.Value = "'1.176"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +0.30%  "
$c = $ws.Range("D37")
$c.Value = "'2.695"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +0.50%  "
$ws.Range("D38").Value = "1.293.88"
$ws.Range("E38").Value = "  +11.90%  "
$c = $ws.Range("D39")
$c.Value = "'0.01883"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +1.33%  "
$c = $ws.Range("D40")
$c.Value = "'2.738"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +0.82%  "
$c = $ws.Range("D41")
$c.Value = "'6.584"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +3.39%  "
$c = $ws.Range("D42")
$c.Value = "'110.98"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +5.06%  "
$c = $ws.Range("D43")
$c.Value = "'0.8960"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +1.20%  "
$ws.Range("E44").Value = "  +0.26%  "
$c = $ws.Range("D45")
$c.Value = "'0.00000000131"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +8.64%  "
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("D47").Value = "2.016.59"
$ws.Range("E47").Value = "  -0.46%  "
$c = $ws.Range("D48")
$c.Value = "'1.802"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -0.10%  "
$c = $ws.Range("D49")
$c.Value = "'0.5202"
$c.Style = "Normal"
$c = $ws.Range("D50")
$c.Value = "'9.449"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +0.71%  "
$c = $ws.Range("D51")
$c.Value = "'0.4358"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +1.31%  "
